# edit.ps1 - applies the "cleaning up the last edits" changes to
# land_classification_diagram.pptx
#
#   1. Refresh the cached "datetimeFigureOut" field text (2/19/2025 -> 2/24/2025)
#      on the slide master and every slide layout's Date placeholder.
#   2. Slide 1, "Proxy Analysis Forest Land Base" box:
#        - title run: drop the leading "Provincial " so it reads
#          "Proxy Analysis Forest Land Base"
#        - description paragraph: append "that are analyzed during
#          timber supply review." (the last two words get their own,
#          newly-italicised run, matching a fresh edit boundary)
#   3. Slide 1, "Proxy Timber Harvesting Land Base" box:
#        - title run: "Provincial" -> "Proxy"
#        - description paragraph: "available" -> "acceptable"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder cache refresh (slide master + all slide layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2/24/2025"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2 & 3. Slide 1 text edits
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)

# --- "Proxy Analysis Forest Land Base" box (2nd shape in the group) ---
$pAFLB = $group.GroupItems.Item(2)
$trAFLB = $pAFLB.TextFrame2.TextRange

# Title: "Provincial Proxy Analysis Forest Land Base" -> "Proxy Analysis Forest Land Base"
$titleLead = $trAFLB.Characters(1, 11)
$titleLead.Text = ""

# Description paragraph: replace the whole run with the extended sentence.
$fullAFLB = $trAFLB.Text
$descStart0 = $fullAFLB.IndexOf("The portion of the FALB")
$descStartIdx = $descStart0 + 1
$descLen = $fullAFLB.Length - $descStart0
$newDescAFLB = "The portion of the FALB that is managed for timber supply and other forest management objectives (e.g., landscape-level biodiversity) that are analyzed during timber supply review."
$descRange = $trAFLB.Characters($descStartIdx, $descLen)
$descRange.Text = $newDescAFLB

# Give "timber supply review." its own run (as happens when that phrase is
# typed in as a later, separate edit) by nudging its formatting.
$fullAFLB2 = $trAFLB.Text
$tailStart0 = $fullAFLB2.IndexOf("timber supply review.")
$tailStartIdx = $tailStart0 + 1
$tailLen = "timber supply review.".Length
$tailRange = $trAFLB.Characters($tailStartIdx, $tailLen)
$tailRange.Font.Italic = 1

# --- "Proxy Timber Harvesting Land Base" box (3rd shape in the group) ---
$pTHLB = $group.GroupItems.Item(3)
$trTHLB = $pTHLB.TextFrame2.TextRange

# Title: "Provincial Timber Harvesting Land Base" -> "Proxy Timber Harvesting Land Base"
$titleTHLB = $trTHLB.Characters(1, 38)
$titleTHLB.Text = "Proxy Timber Harvesting Land Base"

# Description paragraph: "available" -> "acceptable" (replace whole run
# so the run stays a single <a:r> like the source).
$fullTHLB = $trTHLB.Text
$desc2Start0 = $fullTHLB.IndexOf("The portion of the AFLB")
$desc2StartIdx = $desc2Start0 + 1
$desc2Len = $fullTHLB.Length - $desc2Start0
$newDescTHLB = "The portion of the AFLB where timber harvesting is considered both acceptable and economically feasible, given the objectives for all relevant forest values, existing timber quality, market values and applicable technology. It includes areas where timber harvesting is limited due to management objectives (e.g., conditional harvest zones)."
$desc2Range = $trTHLB.Characters($desc2StartIdx, $desc2Len)
$desc2Range.Text = $newDescTHLB
